$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60, shifting existing rows 60..142 down to 61..143
$ws.Rows.Item(60).Insert()

# Populate the new row 60 with data
$ws.Range("A60").Value = 9
$ws.Range("B60").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C60").Value = "Metropolitana"
$ws.Range("D60").Value = 44848
$ws.Range("E60").Value = 13
$ws.Range("F60").Value = 100112022
$ws.Range("G60").Value = "Arveja Verde"
$ws.Range("H60").Value = "Perfection"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 40
$ws.Range("K60").Value = 28000
$ws.Range("L60").Value = 28000
$ws.Range("M60").Value = 28000
$ws.Range("N60").Value = "`$/malla 25 kilos"
$ws.Range("O60").Value = "Provincia de Limarí"
$ws.Range("P60").Value = 1120
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"
